$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "Y"
$ws.Range("D14").Value = "Y"

$ws.Range("D4").HorizontalAlignment = -4108
$ws.Range("D4").VerticalAlignment = -4108
$ws.Range("D14").HorizontalAlignment = -4108
$ws.Range("D14").VerticalAlignment = -4108

$ws.Range("E12").Select()
